# Commit: "added doc until motif analysis"
#
# - Rename the "peals" sheet (sheetId 2) to "peaks" (typo fix).
# - Make the "peaks" sheet the active tab instead of "differential peaks",
#   leaving the selection on cell C28 there.
# - "differential peaks" stops being the tab-selected sheet; its own
#   selection (A2:B3) is left untouched.

$wb = $excel.ActiveWorkbook

$peaksSheet = $wb.Worksheets.Item("peals")
$peaksSheet.Name = "peaks"

$peaksSheet.Activate()
$peaksSheet.Range("C28").Select() | Out-Null
